# Fruta / hortaliza, semanal
# Inserts two new weekly price rows (Limon, Feria Lagunitas de Puerto Montt)
# right above the old row 231, pushing the existing rows 231-273 down to
# 233-275, and populates the two new rows with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 231 (old row 231 data shifts to 233).
$ws.Rows.Item(231).Insert()
$ws.Rows.Item(231).Insert()

# Common columns (A-C, E-K) repeat the same market/category block used by
# every row around them.
$commonCols = @{
    1  = 4                                      # A Mercado ID
    2  = "Feria Lagunitas de Puerto Montt"       # B Mercado
    3  = "Los Lagos"                             # C Region
    5  = 10                                      # E Codreg
    6  = "Fruta"                                 # F Tipo
    7  = 100102                                  # G Producto ID
    8  = "Cítricos"                              # H Producto
    9  = 100102003                               # I Categoria ID
    10 = "Limón"                                 # J Categoria
    11 = "Sin especificar"                       # K Variedad
}

$newRows = @(
    @{ Row = 231; D = 44474; L = "1a amarillo"; M = 800; N = 9000; O = 9000; P = 9000; Q = "`$/malla 18 kilos"; R = "Provincia de Melipilla"; S = 500; T = 18 }
    @{ Row = 232; D = 44474; L = "2a amarillo"; M = 600; N = 7000; O = 7000; P = 7000; Q = "`$/malla 18 kilos"; R = "Provincia de Melipilla"; S = 389; T = 18 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    foreach ($col in $commonCols.Keys) {
        $ws.Cells.Item($r, $col).Value = $commonCols[$col]
    }

    $ws.Cells.Item($r, 4).Value = $rowData.D    # D Fecha
    $ws.Cells.Item($r, 12).Value = $rowData.L   # L Calidad
    $ws.Cells.Item($r, 13).Value = $rowData.M   # M Volumen
    $ws.Cells.Item($r, 14).Value = $rowData.N   # N Precio minimo
    $ws.Cells.Item($r, 15).Value = $rowData.O   # O Precio maximo
    $ws.Cells.Item($r, 16).Value = $rowData.P   # P Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $rowData.Q   # Q Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $rowData.R   # R Origen
    $ws.Cells.Item($r, 19).Value = $rowData.S   # S Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $rowData.T   # T Kg / unidad
}
